$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ------------------------------------------------------------------
# Sheet1 ("杨瀚森")
# ------------------------------------------------------------------

# Row 4: bump review count + roll the three review dates forward
$ws1.Range("D4").Value = 4
$ws1.Range("E4").Value = 43083
$ws1.Range("F4").Value = 43098
$ws1.Range("G4").Value = 43128

# Row 7: bump review count + roll the three review dates forward
$ws1.Range("D7").Value = 3
$ws1.Range("E7").Value = 43083
$ws1.Range("F7").Value = 43090
$ws1.Range("G7").Value = 43105

# Row 9: the next two review dates haven't happened yet - clear them
# but keep the date-formatted style on the cells
$ws1.Range("F9").ClearContents()
$ws1.Range("G9").ClearContents()

# Row 10 (高中单词9) is removed entirely - that word list is done
$ws1.Rows.Item(10).Delete()

# The "bestFit" wide column used to be column E, now it is column F
$ws1.Columns.Item(6).ColumnWidth = 19.86

# Move the active selection to the newly-cleared F9:G9 range
$ws1.Range("F9:G9").Select() | Out-Null

# Restore sheet1's print setup
$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1

# ------------------------------------------------------------------
# Sheet2 ("尹嘉禾")
# ------------------------------------------------------------------

# Row 2 (21天list1): it has now been reviewed once more
$ws2.Range("D2").Value = 1
$ws2.Range("E2").Value = 43083
$ws2.Range("F2").Value = 43085
$ws2.Range("G2").Value = 43089

# Copy row 2's formatting down into the new row 3 before filling values
$ws2.Range("A2:G2").Copy() | Out-Null
$ws2.Range("A3:G3").PasteSpecial(-4122) | Out-Null

# Row 3 (new): 21天list2 just started
$ws2.Range("A3").Value = 1
$ws2.Range("B3").Value = 43083
$ws2.Range("C3").Value = "21天list2"
$ws2.Range("D3").Value = 0
$ws2.Range("E3").Value = 43083
$ws2.Range("F3").Value = 43084
$ws2.Range("G3").Value = 43086

# Reset sheet2's selection back to the default A1 and make sheet1 active again
$ws2.Activate()
$ws2.Range("A1").Select() | Out-Null
$ws1.Activate()
